# Populate Sheet1 with a short list of forum-post bookmarks, each one a
# hyperlinked piece of text, as in the target workbook "links.xlsx".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cell text --------------------------------------------------------
$ws.Range("A1").Value = "AUC, MCC, and logloss"
$ws.Range("A2").Value = "corss validation"
$ws.Range("A3").Value = "production map"
$ws.Range("A4").Value = "how many ones"
$ws.Range("A5").Value = "MCC explained"
$ws.Range("A6").Value = "categorical features"
$ws.Range("A7").Value = "best single model"

# --- hyperlinks ---------------------------------------------------------
# (SubAddress is only used on row 5, which links to a specific forum post)
$ws.Hyperlinks.Add($ws.Range("A1"), "https://www.kaggle.com/forums/f/208/bnp-paribas-cardif-claims-management/t/19240/auc-mcc-and-logloss")
$ws.Hyperlinks.Add($ws.Range("A2"), "https://www.kaggle.com/forums/f/208/bnp-paribas-cardif-claims-management/t/19082/corss-validation")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.kaggle.com/forums/f/208/bnp-paribas-cardif-claims-management/t/19220/production-map")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://www.kaggle.com/forums/f/208/bnp-paribas-cardif-claims-management/t/19096/how-many-ones")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://www.kaggle.com/forums/f/208/bnp-paribas-cardif-claims-management/t/19143/mcc-explained", "post138285")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://www.kaggle.com/forums/f/208/bnp-paribas-cardif-claims-management/t/19075/categorical-features")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://www.kaggle.com/forums/f/208/bnp-paribas-cardif-claims-management/t/19322/best-single-model")

# --- selection matches the saved view in the target file ----------------
$null = $ws.Range("A7").Select()
